$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 24.14924052890208

# Row 3
$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 250.4280181796688

# Row 4
$ws.Range("B4").Value = 0.7287194209349384
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 24.53146345053645
